$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new test-case data
$ws.Range("A2").Value = "iAU_TC_ID_245"
$ws.Range("B2").Value = "@RegressionA Validation of Exam Submit for Approval"
$ws.Range("C2").Value = "failed"

# Remove row 3 entirely (was a duplicate-ish "approve" row)
$ws.Rows("3:3").Delete()
